$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "25.927.62"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -4.78%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.820.69"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -4.27%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.005"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  +0.31%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "281.88"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -8.00%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.5056"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -5.64%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3548"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -6.88%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "45.01"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.83%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.06689"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -8.15%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "20.16"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -8.41%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.8513"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -5.49%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.07880"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -3.84%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "1.809.94"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +61.76%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "5.042"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -5.51%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "87.40"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -8.65%  "

$ws.Range("E17").Value = "  +0.42%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "14.05"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -5.19%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.000008131"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -5.86%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.26%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "25.999.73"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -4.64%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.749"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -5.48%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "10.12"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -5.95%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "6.123"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -6.02%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "141.61"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -5.35%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.151"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -5.93%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "1.681"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -3.14%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "16.98"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -7.49%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "108.82"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -6.73%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "4.308"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -10.38%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "4.238"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -11.42%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.08815"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -4.72%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.04792"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -5.22%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.7421"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -10.50%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.124"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -7.82%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "2.858"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -4.81%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.004"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +0.34%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "3.098"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -6.87%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.427"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -9.24%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.5411"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -5.53%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.01858"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -7.18%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.9853"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -8.41%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "112.37"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -4.09%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "6.225"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -5.22%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "8.207"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -11.77%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.4726"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -4.33%  "

$ws.Range("E47").Value = "  +0.30%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.1376"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -9.49%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "9.265"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -8.55%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "35.72"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -6.97%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.498"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -8.40%  "
